# Update readme.pdf and samples:
# - rename header row to lowercase field names (and a few renames)
# - flesh out the "New York's cool new park in the sky" summary row
#   (icon_color + thumb_url were missing, caption text gets a "More Info" link)
# - drop wrap-text on the two pic_url cells that now show the summary row
# - grow row 2's height to fit the longer caption
# - move the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header renames -------------------------------------------------
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "caption"
$ws.Range("C1").Value = "icon_color"
$ws.Range("D1").Value = "address"
$ws.Range("E1").Value = "city"
$ws.Range("F1").Value = "state"
$ws.Range("G1").Value = "zip"
$ws.Range("H1").Value = "pic_url"
$ws.Range("I1").Value = "thumb_url"

# --- Row 2: fill in the previously-missing fields ---------------------------
$ws.Range("C2").Value = "R"
$ws.Range("B2").Value = 'Hovering above Manhattan''s West Side on a formerly abandoned elevated railroad is an aerial park that has become one of New York City''s top attractions. <i><a href="http://www.thehighline.org" style="color:yellow" target="_blank">More Info</a></i>'
$ws.Range("I2").Value = "http://farm8.staticflickr.com/7062/6855356176_f7f5801fd5_m.jpg"

# --- Formatting: drop wrap text on H2/H3 (keep vertical-top alignment) -----
$ws.Range("H2").WrapText = $false
$ws.Range("H3").WrapText = $false

# --- Row height for the taller summary row ----------------------------------
$ws.Rows.Item(2).RowHeight = 79.2

# --- Selection ---------------------------------------------------------------
$ws.Range("F2").Select()
